$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: Facturas (invoices) history entry
$ws.Range("C27").Value = "SP_FACTURASHIS_SELECT"
$ws.Range("A27").Value = "historicofacturas"
$ws.Range("B27").Value = "FacturashisController"
$ws.Range("H27").Value = "ConsultarFacturasDelDia"
$ws.Range("G27").Value = "dailyoverview-page"

# Row 28: Genericos (generic securities) history entry
$ws.Range("A28").Value = "historicogenericos"
$ws.Range("B28").Value = "GenericoshisController"
$ws.Range("C28").Value = "SP_GENERICOSSHIS_SELECT"
$ws.Range("G28").Value = "dailyoverview-page"
$ws.Range("H28").Value = "ConsultarGenericosDelDia"

# Copy style formatting from the row above (row 26) for the new rows
$ws.Range("A26:C26").Copy()
$ws.Range("A27:C28").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("G26:H26").Copy()
$ws.Range("G27:H28").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("H32").Select()
